$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 50,4
$arr[0,0] = 0.1060158684849739
$arr[0,1] = 0.98294597864151
$arr[0,2] = 0.01189486496150494
$arr[0,3] = 0.9982830882072449
$arr[1,0] = 0.03220825642347336
$arr[1,1] = 0.9959748387336731
$arr[1,2] = 0.01809039153158665
$arr[1,3] = 0.9979015588760376
$arr[2,0] = 0.01974356546998024
$arr[2,1] = 0.9966527819633484
$arr[2,2] = 0.007627138402312994
$arr[2,3] = 0.9987600445747375
$arr[3,0] = 0.01533698569983244
$arr[3,1] = 0.996589183807373
$arr[3,2] = 0.004463640507310629
$arr[3,3] = 0.9989507794380188
$arr[4,0] = 0.0120523264631629
$arr[4,1] = 0.9970552921295166
$arr[4,2] = 0.005144278053194284
$arr[4,3] = 0.9981877207756042
$arr[5,0] = 0.008542587049305439
$arr[5,1] = 0.9981145262718201
$arr[5,2] = 0.002715906128287315
$arr[5,3] = 0.9997138381004333
$arr[6,0] = 0.007502995897084475
$arr[6,1] = 0.9982839822769165
$arr[6,2] = 0.002613786607980728
$arr[6,3] = 0.9997138381004333
$arr[7,0] = 0.007005696184933186
$arr[7,1] = 0.9982839822769165
$arr[7,2] = 0.003068094374611974
$arr[7,3] = 0.9994276762008667
$arr[8,0] = 0.005488919094204903
$arr[8,1] = 0.9987500905990601
$arr[8,2] = 0.004755532369017601
$arr[8,3] = 0.9994276762008667
$arr[9,0] = 0.006934571079909801
$arr[9,1] = 0.9984322786331177
$arr[9,2] = 0.00293750362470746
$arr[9,3] = 0.9997138381004333
$arr[10,0] = 0.007033452391624451
$arr[10,1] = 0.9982839822769165
$arr[10,2] = 0.002698276191949844
$arr[10,3] = 0.9997138381004333
$arr[11,0] = 0.006866606883704662
$arr[11,1] = 0.9986229538917542
$arr[11,2] = 0.002456559799611568
$arr[11,3] = 0.9996184706687927
$arr[12,0] = 0.00552320946007967
$arr[12,1] = 0.9987500905990601
$arr[12,2] = 0.002960914047434926
$arr[12,3] = 0.9995231032371521
$arr[13,0] = 0.006903073750436306
$arr[13,1] = 0.9983475804328918
$arr[13,2] = 0.002061939565464854
$arr[13,3] = 0.9997138381004333
$arr[14,0] = 0.005548911169171333
$arr[14,1] = 0.9986865520477295
$arr[14,2] = 0.002646982436999679
$arr[14,3] = 0.999809205532074
$arr[15,0] = 0.006238820031285286
$arr[15,1] = 0.9984534978866577
$arr[15,2] = 0.003029127838090062
$arr[15,3] = 0.9997138381004333
$arr[16,0] = 0.005873973481357098
$arr[16,1] = 0.9985805749893188
$arr[16,2] = 0.003058843547478318
$arr[16,3] = 0.999809205532074
$arr[17,0] = 0.006283506751060486
$arr[17,1] = 0.9985382556915283
$arr[17,2] = 0.00215220102109015
$arr[17,3] = 0.999809205532074
$arr[18,0] = 0.006131339818239212
$arr[18,1] = 0.9985805749893188
$arr[18,2] = 0.002312135649845004
$arr[18,3] = 0.999809205532074
$arr[19,0] = 0.005830490961670876
$arr[19,1] = 0.9986017942428589
$arr[19,2] = 0.002500922651961446
$arr[19,3] = 0.9997138381004333
$arr[20,0] = 0.005165109876543283
$arr[20,1] = 0.9988348484039307
$arr[20,2] = 0.002457036869600415
$arr[20,3] = 0.999809205532074
$arr[21,0] = 0.005619572475552559
$arr[21,1] = 0.9986865520477295
$arr[21,2] = 0.003119073109701276
$arr[21,3] = 0.999809205532074
$arr[22,0] = 0.004923094529658556
$arr[22,1] = 0.9988136291503906
$arr[22,2] = 0.002870586700737476
$arr[22,3] = 0.999809205532074
$arr[23,0] = 0.006507820449769497
$arr[23,1] = 0.9985805749893188
$arr[23,2] = 0.002109257271513343
$arr[23,3] = 0.999809205532074
$arr[24,0] = 0.00503992848098278
$arr[24,1] = 0.9987712502479553
$arr[24,2] = 0.001816458883695304
$arr[24,3] = 0.999809205532074
$arr[25,0] = 0.005642907693982124
$arr[25,1] = 0.99872887134552
$arr[25,2] = 0.002152129309251904
$arr[25,3] = 0.999809205532074
$arr[26,0] = 0.006288413424044847
$arr[26,1] = 0.9984322786331177
$arr[26,2] = 0.003456659615039825
$arr[26,3] = 0.9997138381004333
$arr[27,0] = 0.005495802499353886
$arr[27,1] = 0.9986653327941895
$arr[27,2] = 0.003089317586272955
$arr[27,3] = 0.999809205532074
$arr[28,0] = 0.004574081394821405
$arr[28,1] = 0.9988983869552612
$arr[28,2] = 0.003325679106637836
$arr[28,3] = 0.999809205532074
$arr[29,0] = 0.005157672334462404
$arr[29,1] = 0.9987924695014954
$arr[29,2] = 0.003631195984780788
$arr[29,3] = 0.999809205532074
$arr[30,0] = 0.006275206338614225
$arr[30,1] = 0.9985170364379883
$arr[30,2] = 0.004379638005048037
$arr[30,3] = 0.9997138381004333
$arr[31,0] = 0.005880849901586771
$arr[31,1] = 0.9985805749893188
$arr[31,2] = 0.003731110598891973
$arr[31,3] = 0.999809205532074
$arr[32,0] = 0.004843940027058125
$arr[32,1] = 0.9988560080528259
$arr[32,2] = 0.003582606092095375
$arr[32,3] = 0.999809205532074
$arr[33,0] = 0.005134676583111286
$arr[33,1] = 0.9987077116966248
$arr[33,2] = 0.003486247500404716
$arr[33,3] = 0.999809205532074
$arr[34,0] = 0.005623715929687023
$arr[34,1] = 0.9986017942428589
$arr[34,2] = 0.004535359796136618
$arr[34,3] = 0.9997138381004333
$arr[35,0] = 0.004520641639828682
$arr[35,1] = 0.9988348484039307
$arr[35,2] = 0.004230548162013292
$arr[35,3] = 0.9997138381004333
$arr[36,0] = 0.004571388941258192
$arr[36,1] = 0.9987500905990601
$arr[36,2] = 0.003789717564359307
$arr[36,3] = 0.999809205532074
$arr[37,0] = 0.005260188598185778
$arr[37,1] = 0.9984111189842224
$arr[37,2] = 0.003331918036565185
$arr[37,3] = 0.999809205532074
$arr[38,0] = 0.004574465099722147
$arr[38,1] = 0.9985382556915283
$arr[38,2] = 0.003652730258181691
$arr[38,3] = 0.9997138381004333
$arr[39,0] = 0.003656366607174277
$arr[39,1] = 0.9989195466041565
$arr[39,2] = 0.004458204377442598
$arr[39,3] = 0.9997138381004333
$arr[40,0] = 0.003555480390787125
$arr[40,1] = 0.9991313815116882
$arr[40,2] = 0.004676146432757378
$arr[40,3] = 0.999809205532074
$arr[41,0] = 0.003170523094013333
$arr[41,1] = 0.9992373585700989
$arr[41,2] = 0.004576669074594975
$arr[41,3] = 0.999809205532074
$arr[42,0] = 0.004053585696965456
$arr[42,1] = 0.9989831447601318
$arr[42,2] = 0.005633137654513121
$arr[42,3] = 0.999809205532074
$arr[43,0] = 0.005070705432444811
$arr[43,1] = 0.9987712502479553
$arr[43,2] = 0.006026525516062975
$arr[43,3] = 0.9997138381004333
$arr[44,0] = 0.003692806698381901
$arr[44,1] = 0.9991313815116882
$arr[44,2] = 0.005791679490357637
$arr[44,3] = 0.999809205532074
$arr[45,0] = 0.003790432587265968
$arr[45,1] = 0.9990466833114624
$arr[45,2] = 0.006063591688871384
$arr[45,3] = 0.999809205532074
$arr[46,0] = 0.003532886737957597
$arr[46,1] = 0.9991949796676636
$arr[46,2] = 0.006573354359716177
$arr[46,3] = 0.999809205532074
$arr[47,0] = 0.003735285717993975
$arr[47,1] = 0.9991949796676636
$arr[47,2] = 0.007104669697582722
$arr[47,3] = 0.9997138381004333
$arr[48,0] = 0.002962330589070916
$arr[48,1] = 0.9993220567703247
$arr[48,2] = 0.007400864735245705
$arr[48,3] = 0.9997138381004333
$arr[49,0] = 0.002940981416031718
$arr[49,1] = 0.9993432760238647
$arr[49,2] = 0.006460781674832106
$arr[49,3] = 0.999809205532074
$ws.Range("A2:D51").Value = $arr
